$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 165.66667
$ws.Range("I6").Value = 165.66667
$ws.Range("K6").Value = 497.00001
$ws.Range("M6").Value = -385.00001

$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws.Range("H17").Value = 958.6667
$ws.Range("J17").Value = 835.625
$ws.Range("L17").Value = 2506.875
$ws.Range("N17").Value = -2842.875

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H111").Value = 1757.8
$ws.Range("I111").Value = 1198
$ws.Range("K111").Value = 3594
$ws.Range("M111").Value = -527

$ws.Range("H112").Value = 5030
$ws.Range("I112").Value = 1150
$ws.Range("J112").Value = 6000
$ws.Range("K112").Value = 3450
$ws.Range("L112").Value = 18000
$ws.Range("M112").Value = -2342
$ws.Range("N112").Value = -20216


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 958.3333
$ws.Range("J29").Value = 2000
$ws.Range("L29").Value = 2000
$ws.Range("N29").Value = -2616

$ws.Range("H32").Value = 3445.7827
$ws.Range("I32").Value = 2112.775
$ws.Range("J32").Value = 12332.5
$ws.Range("K32").Value = 2112.775
$ws.Range("L32").Value = 12332.5
$ws.Range("M32").Value = -1825.775
$ws.Range("N32").Value = -12906.5

$ws.Range("H76").Value = 30499.75
$ws.Range("J76").Value = 30499.75
$ws.Range("L76").Value = 30499.75
$ws.Range("N76").Value = -31175.75

$ws.Range("H79").Value = 30499.75
$ws.Range("J79").Value = 30499.75
$ws.Range("L79").Value = 30499.75
$ws.Range("N79").Value = -32839.75

$ws.Range("H102").Value = 27778236
$ws.Range("J102").Value = 615
$ws.Range("L102").Value = 615
$ws.Range("N102").Value = -3859


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3208695.2
$ws.Range("I105").Value = 4905027
$ws.Range("J105").Value = 4513
$ws.Range("K105").Value = 4905027
$ws.Range("L105").Value = 4513
$ws.Range("M105").Value = -4903280
$ws.Range("N105").Value = -8007

$ws.Range("H134").Value = 3221.25
$ws.Range("I134").Value = 2515.5
$ws.Range("K134").Value = 7546.5
$ws.Range("M134").Value = -5011.5


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 674.875
$ws.Range("I2").Value = 839.8
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 839.8
$ws.Range("L2").Value = 400
$ws.Range("M2").Value = -726.8
$ws.Range("N2").Value = -626

$ws.Range("H22").Value = 67899.95
$ws.Range("I22").Value = 86749.86
$ws.Range("J22").Value = 23916.834
$ws.Range("K22").Value = 86749.86
$ws.Range("L22").Value = 23916.834
$ws.Range("M22").Value = -86399.86
$ws.Range("N22").Value = -24616.834

$ws.Range("H62").Value = 136519.67
$ws.Range("I62").Value = 4780
$ws.Range("K62").Value = 4780
$ws.Range("M62").Value = -4156

$ws.Range("H65").Value = 136519.67
$ws.Range("I65").Value = 4780
$ws.Range("K65").Value = 23900
$ws.Range("M65").Value = -20780

$ws.Range("H86").Value = 5000
$ws.Range("I86").Value = 5000
$ws.Range("K86").Value = 5000
$ws.Range("M86").Value = -3877

$ws.Range("H89").Value = 5000
$ws.Range("I89").Value = 5000
$ws.Range("K89").Value = 25000
$ws.Range("M89").Value = -19384


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 500
$ws.Range("I31").Value = 500
$ws.Range("K31").Value = 1500
$ws.Range("M31").Value = -1212

$ws.Range("H39").Value = 8991.523
$ws.Range("J39").Value = 11021.467
$ws.Range("L39").Value = 33064.401
$ws.Range("N39").Value = -33652.401

$ws.Range("H107").Value = 166987.17
$ws.Range("I107").Value = 356
$ws.Range("J107").Value = 500249.5
$ws.Range("K107").Value = 1068
$ws.Range("L107").Value = 1500748.5
$ws.Range("M107").Value = 852
$ws.Range("N107").Value = -1504588.5

$ws.Range("H113").Value = 833
$ws.Range("J113").Value = 999.5
$ws.Range("L113").Value = 2998.5
$ws.Range("N113").Value = -7338.5

$ws.Range("H129").Value = 1333.75
$ws.Range("I129").Value = 695
$ws.Range("J129").Value = 3250
$ws.Range("K129").Value = 2085
$ws.Range("L129").Value = 9750
$ws.Range("M129").Value = 2915
$ws.Range("N129").Value = -19750


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 620
$ws.Range("I2").Value = 160
$ws.Range("J2").Value = 850
$ws.Range("K2").Value = 160
$ws.Range("L2").Value = 850
$ws.Range("M2").Value = -47
$ws.Range("N2").Value = -1076

$ws.Range("H9").Value = 200
$ws.Range("I9").Value = 200
$ws.Range("K9").Value = 200
$ws.Range("M9").Value = -30

$ws.Range("H18").Value = 2759576.2
$ws.Range("J18").Value = 19150
$ws.Range("L18").Value = 19150
$ws.Range("N18").Value = -19736

$ws.Range("H21").Value = 12506003
$ws.Range("I21").Value = 25000006
$ws.Range("K21").Value = 25000006
$ws.Range("M21").Value = -24999833

$ws.Range("H29").Value = 19664
$ws.Range("I29").Value = 19664
$ws.Range("K29").Value = 19664
$ws.Range("M29").Value = -19374

$ws.Range("H30").Value = 12506003
$ws.Range("I30").Value = 25000006
$ws.Range("K30").Value = 25000006
$ws.Range("M30").Value = -24999901


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 45998.75
$ws.Range("I23").Value = 42998.332
$ws.Range("J23").Value = 55000
$ws.Range("K23").Value = 42998.332
$ws.Range("L23").Value = 55000
$ws.Range("M23").Value = -42768.332
$ws.Range("N23").Value = -55460

$ws.Range("H40").Value = 4006.5264
$ws.Range("I40").Value = 3954.9333
$ws.Range("K40").Value = 3954.9333
$ws.Range("M40").Value = -3818.9333

$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()

$ws.Range("H131").Value = 7222.222
$ws.Range("J131").Value = 7222.222
$ws.Range("L131").Value = 7222.222
$ws.Range("N131").Value = -17302.222


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3998
$ws.Range("I62").Value = 3998
$ws.Range("K62").Value = 3998
$ws.Range("M62").Value = -3374

$ws.Range("H65").Value = 3998
$ws.Range("I65").Value = 3998
$ws.Range("K65").Value = 19990
$ws.Range("M65").Value = -16870

$ws.Range("H81").Value = 4849.7144
$ws.Range("I81").Value = 5574.8335
$ws.Range("J81").Value = 499
$ws.Range("K81").Value = 11149.667
$ws.Range("L81").Value = 998
$ws.Range("M81").Value = -10088.667
$ws.Range("N81").Value = -3120

$ws.Range("H84").Value = 4849.7144
$ws.Range("I84").Value = 5574.8335
$ws.Range("J84").Value = 499
$ws.Range("K84").Value = 55748.335
$ws.Range("L84").Value = 4990
$ws.Range("M84").Value = -50444.335
$ws.Range("N84").Value = -15598

$ws.Range("H96").Value = 1924.25
$ws.Range("I96").Value = 1924.25
$ws.Range("K96").Value = 1924.25
$ws.Range("M96").Value = -551.25

$ws.Range("H104").Value = 23749.75
$ws.Range("J104").Value = 23749.75
$ws.Range("L104").Value = 23749.75
$ws.Range("N104").Value = -30737.75

$ws.Range("H108").Value = 60000
$ws.Range("J108").Value = 60000
$ws.Range("L108").Value = 60000
$ws.Range("N108").Value = -67680

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()

